$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# The handoff has completed translation kickoff for both locales, so the
# shared "Ready for handoff" status is now "In Translation" everywhere it
# appears (Overview summary columns + each locale's Status column).
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"

$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"

# Re-fit the status columns now that the text is shorter than
# "Ready for handoff".
$ws1.Columns.Item(5).ColumnWidth = 12.576851254417766
$ws1.Columns.Item(6).ColumnWidth = 12.576851254417766
$ws2.Columns.Item(3).ColumnWidth = 12.576851254417766
$ws3.Columns.Item(3).ColumnWidth = 12.576851254417766
